# "version 3 of the generated data. Two signals were still missing -
#  added on request from Rein"
#
# The tripleUuid values (column H, tier01TripleInformation/0/tripleUuid)
# on the "Chemicals & Drugs", "Anatomy" and "Disorders" sheets are
# regenerated, and a handful of pathWeight scores (column A) shift
# slightly as a result of the two added signals.

$wb = $excel.ActiveWorkbook

# --- Chemicals & Drugs --------------------------------------------------
$ws = $wb.Worksheets.Item("Chemicals & Drugs")

$tripleUuids = @{
    2  = "76350720"
    3  = "78841912"
    4  = "116859320"
    5  = "74231788"
    6  = "58409705"
    7  = "56773882"
    8  = "75949843"
    9  = "114043056"
    10 = "61382531"
    11 = "94691585"
    12 = "58031591"
    13 = "53339440"
    14 = "125280756"
    15 = "114040142"
}
foreach ($row in $tripleUuids.Keys) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $tripleUuids[$row]
}

$pathWeights = @{
    4  = 5.313835620880127
    8  = 5.013763904571533
    12 = 4.692009449005127
}
foreach ($row in $pathWeights.Keys) {
    $ws.Cells.Item($row, 1).Value = $pathWeights[$row]
}

# --- Anatomy --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Anatomy")

$tripleUuids = @{
    2  = "58031670"
    3  = "59479264"
    4  = "70409366"
    5  = "57099226"
    6  = "59517342"
    7  = "121751908"
    8  = "78783751"
    9  = "87816211"
    10 = "53682816"
    11 = "102897923"
    12 = "108166559"
    13 = "121941152"
    14 = "122848718"
    15 = "130743106"
    16 = "122851169"
}
foreach ($row in $tripleUuids.Keys) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $tripleUuids[$row]
}

# --- Disorders --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Disorders")

$tripleUuids = @{
    2  = "113483617"
    3  = "127281809"
    4  = "57273052"
    5  = "68296895"
    6  = "84843954"
    7  = "55291627"
    8  = "137260534"
    9  = "132797091"
    10 = "66026027"
    11 = "103628962"
    12 = "131050434"
    13 = "68735561"
    14 = "60580724"
    15 = "110982258"
    16 = "88555167"
    17 = "116859123"
    18 = "63359632"
    19 = "138413405"
    20 = "71155385"
}
foreach ($row in $tripleUuids.Keys) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $tripleUuids[$row]
}

$pathWeights = @{
    2 = 5.908326625823975
    5 = 5.6477155685424805
    6 = 5.51275110244751
    7 = 5.342600345611572
}
foreach ($row in $pathWeights.Keys) {
    $ws.Cells.Item($row, 1).Value = $pathWeights[$row]
}
